# pre-close commit - updated flow processing and added transformation groups
# for data processing

$wb = $excel.ActiveWorkbook

$wsGroups = $wb.Worksheets.Item("Groups")
$wsTransformations = $wb.Worksheets.Item("Transformations")

# Update the selection remembered on the Groups sheet
$wsGroups.Activate()
$wsGroups.Range("A2").Select()

# Add the new transformation groups to the Transformations sheet.
# Values are entered in this exact order so that the generated shared-string
# table indices line up with the target workbook.
$wsTransformations.Range("A2").Value = "rem_none_samples"
$wsTransformations.Range("B2").Value = "cat(sg_1_none,sg_2_none,sg_3_none,sg_4_none) .- mean(control_M9_KC)"
$wsTransformations.Range("B3").Value = "cat(sg_1_atc,sg_2_atc,sg_3_atc,sg_4_atc) .- mean(control_M9_KC_atc)"
$wsTransformations.Range("A3").Value = "rem_atc_samples"
$wsTransformations.Range("A4").Value = "rem_IPTG_samples"
$wsTransformations.Range("B4").Value = "cat(sg_1_IPTG,sg_2_IPTG,sg_3_IPTG,sg_4_IPTG) .- mean(control_M9_KC_IPTG)"
$wsTransformations.Range("A5").Value = "rem_atc_IPTG_samples"
$wsTransformations.Range("B5").Value = "cat(sg_1__atc_IPTG,sg_2_atc_IPTG,sg_3_atc_IPTG,sg_4_atc_IPTG) .- mean(control_M9_KC_atc_IPTG)"
$wsTransformations.Range("A6").Value = "norm_flo"

# Widen column A to fit the new, longer labels
$wsTransformations.Columns.Item(1).ColumnWidth = 19

# Make Transformations the active sheet/tab with B2 selected
$wsTransformations.Activate()
$wsTransformations.Range("B2").Select()
